$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 10.06655267360758
$ws.Cells.Item(2, 3).Value2 = 5.270437362646262
$ws.Cells.Item(2, 4).Value2 = 9.193505539851932
$ws.Cells.Item(2, 5).Value2 = 13.69638693630308
$ws.Cells.Item(2, 6).Value2 = 33.17022212974818
$ws.Cells.Item(2, 8).Value2 = 7.344005520526261
$ws.Cells.Item(2, 9).Value2 = 22.91223486194537
$ws.Cells.Item(2, 10).Value2 = 9.952164006458259
$ws.Cells.Item(2, 11).Value2 = 10.24193201029323
$ws.Cells.Item(2, 15).Value2 = 25.08334463714856

$ws.Cells.Item(3, 2).Value2 = 9.754891060528777
$ws.Cells.Item(3, 3).Value2 = 5.048455835676679
$ws.Cells.Item(3, 4).Value2 = 9.122094076288214
$ws.Cells.Item(3, 5).Value2 = 13.62785279992243
$ws.Cells.Item(3, 6).Value2 = 33.24478385495913
$ws.Cells.Item(3, 8).Value2 = 7.344005520526261
$ws.Cells.Item(3, 9).Value2 = 23.01950234982142
$ws.Cells.Item(3, 10).Value2 = 9.959068111991913
$ws.Cells.Item(3, 11).Value2 = 10.02936680492359
$ws.Cells.Item(3, 15).Value2 = 25.17842383260056

$ws.Cells.Item(4, 2).Value2 = 9.559559341184421
$ws.Cells.Item(4, 3).Value2 = 4.90745107283668
$ws.Cells.Item(4, 4).Value2 = 9.079538283194486
$ws.Cells.Item(4, 5).Value2 = 13.58844026337651
$ws.Cells.Item(4, 6).Value2 = 33.29856291471738
$ws.Cells.Item(4, 8).Value2 = 7.344005520526261
$ws.Cells.Item(4, 9).Value2 = 23.08975856074216
$ws.Cells.Item(4, 10).Value2 = 9.964876781477788
$ws.Cells.Item(4, 11).Value2 = 9.897778754030908
$ws.Cells.Item(4, 15).Value2 = 25.2420919059431

$ws.Cells.Item(5, 2).Value2 = 9.479081855602804
$ws.Cells.Item(5, 3).Value2 = 4.848885965622329
$ws.Cells.Item(5, 4).Value2 = 9.062536574309602
$ws.Cells.Item(5, 5).Value2 = 13.57306317562707
$ws.Cells.Item(5, 6).Value2 = 33.32248484553276
$ws.Cells.Item(5, 8).Value2 = 7.344005520526261
$ws.Cells.Item(5, 9).Value2 = 23.11949309941544
$ws.Cells.Item(5, 10).Value2 = 9.967638875806639
$ws.Cells.Item(5, 11).Value2 = 9.843956709547026
$ws.Cells.Item(5, 15).Value2 = 25.26936480419921

$ws.Cells.Item(6, 2).Value2 = 9.465669276293486
$ws.Cells.Item(6, 3).Value2 = 4.839097078100608
$ws.Cells.Item(6, 4).Value2 = 9.059734436339079
$ws.Cells.Item(6, 5).Value2 = 13.57055148619587
$ws.Cells.Item(6, 6).Value2 = 33.32657810903486
$ws.Cells.Item(6, 8).Value2 = 7.344005520526261
$ws.Cells.Item(6, 9).Value2 = 23.12449720394243
$ws.Cells.Item(6, 10).Value2 = 9.96812138542761
$ws.Cells.Item(6, 11).Value2 = 9.835009701338864
$ws.Cells.Item(6, 15).Value2 = 25.27397356552115

$ws.Cells.Item(7, 2).Value2 = 9.558477384161737
$ws.Cells.Item(7, 3).Value2 = 4.906665603393836
$ws.Cells.Item(7, 4).Value2 = 9.079307595065297
$ws.Cells.Item(7, 5).Value2 = 13.58823009723874
$ws.Cells.Item(7, 6).Value2 = 33.29887741635697
$ws.Cells.Item(7, 8).Value2 = 7.344005520526261
$ws.Cells.Item(7, 9).Value2 = 23.09015509908144
$ws.Cells.Item(7, 10).Value2 = 9.964912432268781
$ws.Cells.Item(7, 11).Value2 = 9.89705360111598
$ws.Cells.Item(7, 15).Value2 = 25.24245434534609

$ws.Cells.Item(8, 2).Value2 = 9.959985004558247
$ws.Cells.Item(8, 3).Value2 = 5.19491949056357
$ws.Cells.Item(8, 4).Value2 = 9.168624334731103
$ws.Cells.Item(8, 5).Value2 = 13.67221078280905
$ws.Cells.Item(8, 6).Value2 = 33.19426843697972
$ws.Cells.Item(8, 8).Value2 = 7.344005520526261
$ws.Cells.Item(8, 9).Value2 = 22.94830866967538
$ws.Cells.Item(8, 10).Value2 = 9.954219023890143
$ws.Cells.Item(8, 11).Value2 = 10.16890378237207
$ws.Cells.Item(8, 15).Value2 = 25.11502880654696

$ws.Cells.Item(9, 2).Value2 = 10.71102946378259
$ws.Cells.Item(9, 3).Value2 = 5.71985281713954
$ws.Cells.Item(9, 4).Value2 = 9.35331483096676
$ws.Cells.Item(9, 5).Value2 = 13.85744888539973
$ws.Cells.Item(9, 6).Value2 = 33.05277596353512
$ws.Cells.Item(9, 8).Value2 = 7.344005520526261
$ws.Cells.Item(9, 9).Value2 = 22.7050192566066
$ws.Cells.Item(9, 10).Value2 = 9.945684373866907
$ws.Cells.Item(9, 11).Value2 = 10.69051175024294
$ws.Cells.Item(9, 15).Value2 = 24.90721715727799

$ws.Cells.Item(10, 2).Value2 = 11.23490512677165
$ws.Cells.Item(10, 3).Value2 = 6.077566056409397
$ws.Cells.Item(10, 4).Value2 = 9.493850490543325
$ws.Cells.Item(10, 5).Value2 = 14.00519182248768
$ws.Cells.Item(10, 6).Value2 = 32.98784597866941
$ws.Cells.Item(10, 8).Value2 = 7.344005520526261
$ws.Cells.Item(10, 9).Value2 = 22.5475469721805
$ws.Cells.Item(10, 10).Value2 = 9.9469653252693
$ws.Cells.Item(10, 11).Value2 = 11.06294533393313
$ws.Cells.Item(10, 15).Value2 = 24.78032356280901

$ws.Cells.Item(11, 2).Value2 = 11.46612981514399
$ws.Cells.Item(11, 3).Value2 = 6.23367511151528
$ws.Cells.Item(11, 4).Value2 = 9.558621904290014
$ws.Cells.Item(11, 5).Value2 = 14.07472974900863
$ws.Cells.Item(11, 6).Value2 = 32.96681784052326
$ws.Cells.Item(11, 8).Value2 = 7.344005520526261
$ws.Cells.Item(11, 9).Value2 = 22.48053195891932
$ws.Cells.Item(11, 10).Value2 = 9.949179269264361
$ws.Cells.Item(11, 11).Value2 = 11.22929621004412
$ws.Cells.Item(11, 15).Value2 = 24.72822453086924

$ws.Cells.Item(12, 2).Value2 = 11.5525944584967
$ws.Cells.Item(12, 3).Value2 = 6.291800860592402
$ws.Cells.Item(12, 4).Value2 = 9.583251850736747
$ws.Cells.Item(12, 5).Value2 = 14.10137870078064
$ws.Cells.Item(12, 6).Value2 = 32.96008050384062
$ws.Cells.Item(12, 8).Value2 = 7.344005520526261
$ws.Cells.Item(12, 9).Value2 = 22.45582003356613
$ws.Cells.Item(12, 10).Value2 = 9.950251258860973
$ws.Cells.Item(12, 11).Value2 = 11.29179090915403
$ws.Cells.Item(12, 15).Value2 = 24.70930729869433

$ws.Cells.Item(13, 2).Value2 = 11.53402257552918
$ws.Cells.Item(13, 3).Value2 = 6.27932698512241
$ws.Cells.Item(13, 4).Value2 = 9.577943091524189
$ws.Cells.Item(13, 5).Value2 = 14.09562559754217
$ws.Cells.Item(13, 6).Value2 = 32.96147697770386
$ws.Cells.Item(13, 8).Value2 = 7.344005520526261
$ws.Cells.Item(13, 9).Value2 = 22.46111258924813
$ws.Cells.Item(13, 10).Value2 = 9.950010012993699
$ws.Cells.Item(13, 11).Value2 = 11.27835460998392
$ws.Cells.Item(13, 15).Value2 = 24.7133453288598

$ws.Cells.Item(14, 2).Value2 = 11.47326566605264
$ws.Cells.Item(14, 3).Value2 = 6.238477176776802
$ws.Cells.Item(14, 4).Value2 = 9.560646278498673
$ws.Cells.Item(14, 5).Value2 = 14.07691593980333
$ws.Cells.Item(14, 6).Value2 = 32.96623898453694
$ws.Cells.Item(14, 8).Value2 = 7.344005520526261
$ws.Cells.Item(14, 9).Value2 = 22.47848555707295
$ws.Cells.Item(14, 10).Value2 = 9.949262786345065
$ws.Cells.Item(14, 11).Value2 = 11.23444799646943
$ws.Cells.Item(14, 15).Value2 = 24.72665191962943

$ws.Cells.Item(15, 2).Value2 = 11.43590557921996
$ws.Cells.Item(15, 3).Value2 = 6.213325619437503
$ws.Cells.Item(15, 4).Value2 = 9.550064268146594
$ws.Cells.Item(15, 5).Value2 = 14.0654963883063
$ws.Cells.Item(15, 6).Value2 = 32.96931550392762
$ws.Cells.Item(15, 8).Value2 = 7.344005520526261
$ws.Cells.Item(15, 9).Value2 = 22.48921365759119
$ws.Cells.Item(15, 10).Value2 = 9.948835481953559
$ws.Cells.Item(15, 11).Value2 = 11.20748734218701
$ws.Cells.Item(15, 15).Value2 = 24.73490835339014

$ws.Cells.Item(16, 2).Value2 = 11.21964473736893
$ws.Cells.Item(16, 3).Value2 = 6.067227445526068
$ws.Cells.Item(16, 4).Value2 = 9.489632868492686
$ws.Cells.Item(16, 5).Value2 = 14.0006926550553
$ws.Cells.Item(16, 6).Value2 = 32.98939162259336
$ws.Cells.Item(16, 8).Value2 = 7.344005520526261
$ws.Cells.Item(16, 9).Value2 = 22.55201959834846
$ws.Cells.Item(16, 10).Value2 = 9.94685338476561
$ws.Cells.Item(16, 11).Value2 = 11.05200714348522
$ws.Cells.Item(16, 15).Value2 = 24.78384178015087

$ws.Cells.Item(17, 2).Value2 = 11.08510351968689
$ws.Cells.Item(17, 3).Value2 = 5.975877426722946
$ws.Cells.Item(17, 4).Value2 = 9.452762420166385
$ws.Cells.Item(17, 5).Value2 = 13.96152121104679
$ws.Cells.Item(17, 6).Value2 = 33.00388861992102
$ws.Cells.Item(17, 8).Value2 = 7.344005520526261
$ws.Cells.Item(17, 9).Value2 = 22.59173288036369
$ws.Cells.Item(17, 10).Value2 = 9.946054591146575
$ws.Cells.Item(17, 11).Value2 = 10.95579664484885
$ws.Cells.Item(17, 15).Value2 = 24.81530354679398

$ws.Cells.Item(18, 2).Value2 = 11.00705570513319
$ws.Cells.Item(18, 3).Value2 = 5.922714400962869
$ws.Cells.Item(18, 4).Value2 = 9.431635791944146
$ws.Cells.Item(18, 5).Value2 = 13.93921118410238
$ws.Cells.Item(18, 6).Value2 = 33.01302772246023
$ws.Cells.Item(18, 8).Value2 = 7.344005520526261
$ws.Cells.Item(18, 9).Value2 = 22.61500970308259
$ws.Cells.Item(18, 10).Value2 = 9.94574875017927
$ws.Cells.Item(18, 11).Value2 = 10.90017303047828
$ws.Cells.Item(18, 15).Value2 = 24.83392886955228

$ws.Cells.Item(19, 2).Value2 = 10.98051858394347
$ws.Cells.Item(19, 3).Value2 = 5.904608924610306
$ws.Cells.Item(19, 4).Value2 = 9.42449702960071
$ws.Cells.Item(19, 5).Value2 = 13.9316957908078
$ws.Cells.Item(19, 6).Value2 = 33.01625953431844
$ws.Cells.Item(19, 8).Value2 = 7.344005520526261
$ws.Cells.Item(19, 9).Value2 = 22.62296548596984
$ws.Cells.Item(19, 10).Value2 = 9.945671606346032
$ws.Cells.Item(19, 11).Value2 = 10.88129253383486
$ws.Cells.Item(19, 15).Value2 = 24.84032591946333

$ws.Cells.Item(20, 2).Value2 = 11.09949488225773
$ws.Cells.Item(20, 3).Value2 = 5.985666347416675
$ws.Cells.Item(20, 4).Value2 = 9.456679167766913
$ws.Cells.Item(20, 5).Value2 = 13.96566840445379
$ws.Cells.Item(20, 6).Value2 = 33.0022624920574
$ws.Cells.Item(20, 8).Value2 = 7.344005520526261
$ws.Cells.Item(20, 9).Value2 = 22.58746032809987
$ws.Cells.Item(20, 10).Value2 = 9.946123733168118
$ws.Cells.Item(20, 11).Value2 = 10.96606843026927
$ws.Cells.Item(20, 15).Value2 = 24.81189958358443

$ws.Cells.Item(21, 2).Value2 = 11.49114172752392
$ws.Cells.Item(21, 3).Value2 = 6.250502880188249
$ws.Cells.Item(21, 4).Value2 = 9.565724140072327
$ws.Cells.Item(21, 5).Value2 = 14.08240298119606
$ws.Cells.Item(21, 6).Value2 = 32.96480699356679
$ws.Cells.Item(21, 8).Value2 = 7.344005520526261
$ws.Cells.Item(21, 9).Value2 = 22.47336463462053
$ws.Cells.Item(21, 10).Value2 = 9.949475932815336
$ws.Cells.Item(21, 11).Value2 = 11.24735840623868
$ws.Cells.Item(21, 15).Value2 = 24.72272140595803

$ws.Cells.Item(22, 2).Value2 = 11.7406928870684
$ws.Cells.Item(22, 3).Value2 = 6.417808458406126
$ws.Cells.Item(22, 4).Value2 = 9.637579688883022
$ws.Cells.Item(22, 5).Value2 = 14.16053135803728
$ws.Cells.Item(22, 6).Value2 = 32.94747164076799
$ws.Cells.Item(22, 8).Value2 = 7.344005520526261
$ws.Cells.Item(22, 9).Value2 = 22.40267427271812
$ws.Cells.Item(22, 10).Value2 = 9.953027960627022
$ws.Cells.Item(22, 11).Value2 = 11.42826664389222
$ws.Cells.Item(22, 15).Value2 = 24.66916975431584

$ws.Cells.Item(23, 2).Value2 = 11.60811242429917
$ws.Cells.Item(23, 3).Value2 = 6.329054291990052
$ws.Cells.Item(23, 4).Value2 = 9.599181325311305
$ws.Cells.Item(23, 5).Value2 = 14.11867104681882
$ws.Cells.Item(23, 6).Value2 = 32.95606967884896
$ws.Cells.Item(23, 8).Value2 = 7.344005520526261
$ws.Cells.Item(23, 9).Value2 = 22.44004793243658
$ws.Cells.Item(23, 10).Value2 = 9.951007975181922
$ws.Cells.Item(23, 11).Value2 = 11.33199860085747
$ws.Cells.Item(23, 15).Value2 = 24.69731750783398

$ws.Cells.Item(24, 2).Value2 = 11.09299071637664
$ws.Cells.Item(24, 3).Value2 = 5.981242781484806
$ws.Cells.Item(24, 4).Value2 = 9.454908184606589
$ws.Cells.Item(24, 5).Value2 = 13.96379280231369
$ws.Cells.Item(24, 6).Value2 = 33.00299515837682
$ws.Cells.Item(24, 8).Value2 = 7.344005520526261
$ws.Cells.Item(24, 9).Value2 = 22.58939056247445
$ws.Cells.Item(24, 10).Value2 = 9.946091996177477
$ws.Cells.Item(24, 11).Value2 = 10.9614255201303
$ws.Cells.Item(24, 15).Value2 = 24.81343684114044

$ws.Cells.Item(25, 2).Value2 = 10.51237233999557
$ws.Cells.Item(25, 3).Value2 = 5.582541935349769
$ws.Cells.Item(25, 4).Value2 = 9.302428742950744
$ws.Cells.Item(25, 5).Value2 = 13.80522453907278
$ws.Cells.Item(25, 6).Value2 = 33.08421428527542
$ws.Cells.Item(25, 8).Value2 = 7.344005520526261
$ws.Cells.Item(25, 9).Value2 = 22.76710104908961
$ws.Cells.Item(25, 10).Value2 = 9.946664542727737
$ws.Cells.Item(25, 11).Value2 = 10.55103660575691
$ws.Cells.Item(25, 15).Value2 = 24.95891737582174
